$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 32.47042033333333
$ws.Range("H2").Value = 97.411261
$ws.Range("I2").Value = 0.5240295449207956
$ws.Range("J2").Value = 0.5240295449207955
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 1106.291454513708
$ws.Range("R2").Value = 9956.623090623378
$ws.Range("S2").Value = 0.2710000185963509
$ws.Range("T2").Value = 0.2710000185963509

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 32.47042033333333
$ws.Range("H3").Value = 97.411261
$ws.Range("I3").Value = 0.5240295449207956
$ws.Range("J3").Value = 0.5240295449207955
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.05649099999999
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 898.9597247272388
$ws.Range("R3").Value = 8090.63752254515
$ws.Range("S3").Value = 0.22021150134034
$ws.Range("T3").Value = 0.22021150134034

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 32.47042033333333
$ws.Range("H4").Value = 97.411261
$ws.Range("I4").Value = 0.5240295449207956
$ws.Range("J4").Value = 0.5240295449207955
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 133.9715797142059
$ws.Range("R4").Value = 1205.744217427853
$ws.Range("S4").Value = 0.0328180249841046
$ws.Range("T4").Value = 0.03281802498410459

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 22.43791
$ws.Range("H5").Value = 67.31372999999999
$ws.Range("I5").Value = 0.3621181261458191
$ws.Range("J5").Value = 0.362118126145819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 764.4763398601631
$ws.Range("R5").Value = 6880.287058741468
$ws.Range("S5").Value = 0.187268103241059
$ws.Range("T5").Value = 0.187268103241059

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 22.43791
$ws.Range("H6").Value = 67.31372999999999
$ws.Range("I6").Value = 0.3621181261458191
$ws.Range("J6").Value = 0.362118126145819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.05649099999999
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 621.2046899912699
$ws.Range("R6").Value = 5590.842209921429
$ws.Range("S6").Value = 0.1521719089964176
$ws.Range("T6").Value = 0.1521719089964176

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 22.43791
$ws.Range("H7").Value = 67.31372999999999
$ws.Range("I7").Value = 0.3621181261458191
$ws.Range("J7").Value = 0.362118126145819
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 92.57786678847665
$ws.Range("R7").Value = 833.2008010962899
$ws.Range("S7").Value = 0.0226781139083424
$ws.Range("T7").Value = 0.02267811390834239

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.054627
$ws.Range("H8").Value = 21.163881
$ws.Range("I8").Value = 0.1138523289333856
$ws.Range("J8").Value = 0.1138523289333855
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 240.3564069932843
$ws.Range("R8").Value = 2163.207662939559
$ws.Range("S8").Value = 0.05887832767682739
$ws.Range("T8").Value = 0.0588783276768274

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.054627
$ws.Range("H9").Value = 21.163881
$ws.Range("I9").Value = 0.1138523289333856
$ws.Range("J9").Value = 0.1138523289333855
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.05649099999999
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 195.310854644619
$ws.Range("R9").Value = 1757.797691801571
$ws.Range("S9").Value = 0.04784385256236747
$ws.Range("T9").Value = 0.04784385256236747

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.054627
$ws.Range("H10").Value = 21.163881
$ws.Range("I10").Value = 0.1138523289333856
$ws.Range("J10").Value = 0.1138523289333855
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 29.10709235612367
$ws.Range("R10").Value = 261.963831205113
$ws.Range("S10").Value = 0.007130148694190671
$ws.Range("T10").Value = 0.00713014869419067
